# B6-PowerPoint.pptx edit:
#   1. Re-style the three tables (slides 14, 15, 16) with the new
#      table-style GUID.
#   2. Swap the presentation's theme colour scheme ("Integral" / Red
#      Violet) for the plain "Office" palette that used to live only in
#      the (otherwise unused) second theme part - this is what the
#      underlying OOXML diff amounts to, since the two theme parts are
#      byte-identical apart from <a:clrScheme> (name + 12 colours) and
#      the top-level theme name.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{86C38E89-9F87-413A-BB9A-2CB71FDB2784}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colour scheme -------------------------------------------
# Colour-scheme slots are addressed 1..12 in the fixed COM order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB longs below are the VBA-style 0x00BBGGRR encoding of the target
# "Office" palette's RRGGBB hex values.
$officeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($slot in $officeColors.Keys) {
    $themeColors.Item($slot).RGB = $officeColors[$slot]
}
